# Refreshes the cryptocurrency price / 1h-volume snapshot (and, for a few
# rows, the coin/link that landed in that rank slot) with the latest scrape
# pulled in by the scheduled GitHub Actions run.
#
# Target cells hold plain text (e.g. "218.71", "26.162.66", "  -0.08%  ")
# even though several values look numeric. Writing a numeric-looking string
# straight into a General-formatted cell would make Excel silently coerce it
# to a floating point number (losing the fixed decimal places / becoming a
# "real" number instead of text), so each cell is switched to Text format
# first and then reset to the "Normal" style afterwards - that keeps the
# stored value a string without leaving behind any stray number formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "26.162.66"
Set-TextCell "E2" "  -0.08%  "

# Row 3
Set-TextCell "D3" "1.656.83"

# Row 4
Set-TextCell "E4" "  -0.30%  "

# Row 5
Set-TextCell "D5" "218.71"
Set-TextCell "E5" "  -0.25%  "

# Row 6
Set-TextCell "D6" "0.5240"
Set-TextCell "E6" "  +0.17%  "

# Row 7
Set-TextCell "D7" "1.003"
Set-TextCell "E7" "  -0.27%  "

# Row 8
Set-TextCell "D8" "0.2666"
Set-TextCell "E8" "  +1.41%  "

# Row 9
Set-TextCell "D9" "0.06348"
Set-TextCell "E9" "  +0.81%  "

# Row 10
Set-TextCell "D10" "20.56"
Set-TextCell "E10" "  -0.12%  "

# Row 11
Set-TextCell "D11" "0.07677"
Set-TextCell "E11" "  -1.82%  "

# Row 12
Set-TextCell "D12" "4.624"
Set-TextCell "E12" "  +2.95%  "

# Row 13
Set-TextCell "D13" "1.653.81"
Set-TextCell "E13" "  -0.34%  "

# Row 14
Set-TextCell "D14" "1.885.57"
Set-TextCell "E14" "  -0.08%  "

# Row 15
Set-TextCell "D15" "0.5618"
Set-TextCell "E15" "  +1.34%  "

# Row 16
Set-TextCell "D16" "0.0₅8205"
Set-TextCell "E16" "  +2.25%  "

# Row 17
Set-TextCell "D17" "65.53"
Set-TextCell "E17" "  +0.58%  "

# Row 18
Set-TextCell "D18" "26.156.08"
Set-TextCell "E18" "  -0.15%  "

# Row 19
Set-TextCell "E19" "  -0.25%  "

# Row 20
Set-TextCell "D20" "4.658"
Set-TextCell "E20" "  +0.53%  "

# Row 21
Set-TextCell "D21" "10.48"
Set-TextCell "E21" "  +3.54%  "

# Row 22
Set-TextCell "D22" "192.31"
Set-TextCell "E22" "  -2.00%  "

# Row 23
Set-TextCell "D23" "5.961"
Set-TextCell "E23" "  +0.19%  "

# Row 24
Set-TextCell "E24" "  -0.33%  "

# Row 25
Set-TextCell "D25" "145.37"
Set-TextCell "E25" "  -0.18%  "

# Row 26
Set-TextCell "D26" "0.1196"
Set-TextCell "E26" "  -0.29%  "

# Row 27
Set-TextCell "D27" "7.275"
Set-TextCell "E27" "  +1.86%  "

# Row 28
Set-TextCell "D28" "15.95"
Set-TextCell "E28" "  -0.45%  "

# Row 29
Set-TextCell "D29" "1.520"
Set-TextCell "E29" "  +1.67%  "

# Row 30
Set-TextCell "D30" "0.05493"
Set-TextCell "E30" "  -4.37%  "

# Row 31
Set-TextCell "E31" "  -0.31%  "

# Row 32
Set-TextCell "D32" "3.474"
Set-TextCell "E32" "  -0.50%  "

# Row 33
Set-TextCell "D33" "3.367"
Set-TextCell "E33" "  -0.05%  "

# Row 34
Set-TextCell "D34" "1.567"
Set-TextCell "E34" "  -1.07%  "

# Row 35
Set-TextCell "D35" "0.9512"
Set-TextCell "E35" "  -0.39%  "

# Row 36
Set-TextCell "E36" "  -0.97%  "

# Row 38
Set-TextCell "D38" "0.5701"
Set-TextCell "E38" "  -0.33%  "

# Row 39
Set-TextCell "D39" "0.01590"
Set-TextCell "E39" "  -0.40%  "

# Row 40
Set-TextCell "D40" "5.881"
Set-TextCell "E40" "  -1.68%  "

# Row 41
Set-TextCell "E41" "  -0.23%  "

# Row 42
Set-TextCell "B42" "TrustWalletToken"
Set-TextCell "C42" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D42" "0.8340"
Set-TextCell "E42" "  -1.60%  "

# Row 43
Set-TextCell "B43" "Maker"
Set-TextCell "C43" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D43" "1.031.15"
Set-TextCell "E43" "  -2.79%  "

# Row 44
Set-TextCell "D44" "101.00"
Set-TextCell "E44" "  -2.84%  "

# Row 45
Set-TextCell "D45" "1.795.55"
Set-TextCell "E45" "  -0.12%  "

# Row 46
Set-TextCell "D46" "58.23"
Set-TextCell "E46" "  +0.07%  "

# Row 47
Set-TextCell "B47" "BabyDogeCoin"
Set-TextCell "C47" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D47" "0.0₈104"
Set-TextCell "E47" "  -1.82%  "

# Row 48
Set-TextCell "B48" "Frax"
Set-TextCell "C48" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell "D48" "0.9999"
Set-TextCell "E48" "  -0.83%  "

# Row 49
Set-TextCell "B49" "EnergySwap"
Set-TextCell "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D49" "8.026"
Set-TextCell "E49" "  +0.27%  "

# Row 50
Set-TextCell "B50" "Mantle"
Set-TextCell "C50" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D50" "0.4345"
Set-TextCell "E50" "  -1.42%  "

# Row 51
Set-TextCell "B51" "Cronos"
Set-TextCell "C51" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D51" "0.05225"
Set-TextCell "E51" "  +0.42%  "
